$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.043.81"
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = "'1.830.38"
$ws.Range('D4').Value = "'0.9987"
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'241.30"
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').Value = "'0.6265"
$ws.Range('E6').Value = '  -5.05%  '
$ws.Range('D7').Value = "'1.000"
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = "'0.07619"
$ws.Range('E8').Value = '  +2.72%  '
$ws.Range('D9').Value = "'45.17"
$ws.Range('E9').Value = '  +8.01%  '
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('D11').Value = "'22.76"
$ws.Range('E11').Value = '  -0.58%  '
$ws.Range('D12').Value = "'0.07639"
$ws.Range('E12').Value = '  -1.67%  '
$ws.Range('D13').Value = "'1.828.92"
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('D14').Value = "'4.955"
$ws.Range('E14').Value = '  -0.48%  '
$ws.Range('D15').Value = "'0.6650"
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('D16').Value = "'82.28"
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('D17').Value = "'0.000009158"
$ws.Range('E17').Value = '  +7.37%  '
$ws.Range('D18').Value = "'5.986"
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').Value = "'28.866.58"
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').Value = "'224.62"
$ws.Range('E20').Value = '  -0.91%  '
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('D22').Value = "'0.9999"
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').Value = "'7.195"
$ws.Range('E23').Value = '  +1.64%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = "'159.83"
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = "'8.418"
$ws.Range('E26').Value = '  -2.06%  '
$ws.Range('E27').Value = '  -2.47%  '
$ws.Range('D28').Value = "'17.80"
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('D29').Value = "'1.494"
$ws.Range('E29').Value = '  -1.45%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = "'4.049"
$ws.Range('E30').Value = '  -1.49%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = "'4.029"
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('D32').Value = "'1.205"
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('D33').Value = "'0.05197"
$ws.Range('E33').Value = '  -1.27%  '
$ws.Range('D34').Value = "'1.847"
$ws.Range('E34').Value = '  -0.94%  '
$ws.Range('E35').Value = '  +0.83%  '
$ws.Range('D36').Value = "'0.7318"
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('D37').Value = "'2.616"
$ws.Range('E37').Value = '  -1.50%  '
$ws.Range('D38').Value = "'1.278.28"
$ws.Range('E38').Value = '  -1.67%  '
$ws.Range('D39').Value = "'2.760"
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('D41').Value = "'6.532"
$ws.Range('E41').Value = '  +8.14%  '
$ws.Range('D42').Value = "'0.8889"
$ws.Range('E42').Value = '  -3.12%  '
$ws.Range('D43').Value = "'0.9998"
$ws.Range('D44').Value = "'101.51"
$ws.Range('E44').Value = '  -1.19%  '
$ws.Range('D45').Value = "'1.976.16"
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('E46').Value = '  -0.67%  '
$ws.Range('D47').Value = "'63.66"
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('D49').Value = "'0.3980"
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('D50').Value = "'0.07294"
$ws.Range('E50').Value = '  -13.75%  '
$ws.Range('D51').Value = "'8.819"
$ws.Range('E51').Value = '  +1.55%  '
